$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C13").Value = 10/1440
$ws.Range("C14").Select()
